$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column header ---
$ws.Range("B1").Value = "Folder Path"

# --- Rows 3-6: column C filled first with the quoted concatenation text ---
$ws.Range("C3").Value = '"Website ID 56"+"\"+"SSA-827\"'
$ws.Range("C4").Value = '"Website ID 56"+"\"+"SSA-1696\"'
$ws.Range("C5").Value = '"Website ID 56"+"\"+"SSA-8000\"'
$ws.Range("C6").Value = '"Website ID 56"+"\"+"Fillable SSA-8000\"'

# --- Rows 7-16: column B then column C (same value) ---
$ws.Range("B7").Value = '"Website ID 56"+"\"+"SSA-8001\"'
$ws.Range("C7").Value = '"Website ID 56"+"\"+"SSA-8001\"'

$ws.Range("B8").Value = '"Website ID 56"+"\"+"SSA-3369\"'
$ws.Range("C8").Value = '"Website ID 56"+"\"+"SSA-3369\"'

$ws.Range("B9").Value = '"Website ID 56"+"\"+"SSA-821\"'
$ws.Range("C9").Value = '"Website ID 56"+"\"+"SSA-821\"'

$ws.Range("B10").Value = '"Website ID 56"+"\"+"HA-1151\"'
$ws.Range("C10").Value = '"Website ID 56"+"\"+"HA-1151\"'

$ws.Range("B11").Value = '"Website ID 56"+"\"+"HA-1152\"'
$ws.Range("C11").Value = '"Website ID 56"+"\"+"HA-1152\"'

$ws.Range("B12").Value = '"Website ID 56"+"\"+"SSA-787\"'
$ws.Range("C12").Value = '"Website ID 56"+"\"+"SSA-787\"'

$ws.Range("B13").Value = '"Website ID 56"+"\"+"SSA-1699\"'
$ws.Range("C13").Value = '"Website ID 56"+"\"+"SSA-1699\"'

$ws.Range("B14").Value = '"Website ID 56"+"\"+"paper form\"'
$ws.Range("C14").Value = '"Website ID 56"+"\"+"paper form\"'

$ws.Range("B15").Value = '"Website ID 56"+"\"+"paper form\"'
$ws.Range("C15").Value = '"Website ID 56"+"\"+"paper form\"'

$ws.Range("B16").Value = '"Website ID 56"+"\"+"SSA-3373\"'
$ws.Range("C16").Value = '"Website ID 56"+"\"+"SSA-3373\"'

# --- Rows 20-26, 29: column B then column C (same value) ---
$ws.Range("B20").Value = '"Website ID 56"+"\"+"SSA-3820\"'
$ws.Range("C20").Value = '"Website ID 56"+"\"+"SSA-3820\"'

$ws.Range("B21").Value = '"Website ID 56"+"\"+"SSA-3375\"'
$ws.Range("C21").Value = '"Website ID 56"+"\"+"SSA-3375\"'

$ws.Range("B22").Value = '"Website ID 56"+"\"+"3376\"'
$ws.Range("C22").Value = '"Website ID 56"+"\"+"3376\"'

$ws.Range("B23").Value = '"Website ID 56"+"\"+"3377\"'
$ws.Range("C23").Value = '"Website ID 56"+"\"+"3377\"'

$ws.Range("B24").Value = '"Website ID 56"+"\"+"3378\"'
$ws.Range("C24").Value = '"Website ID 56"+"\"+"3378\"'

$ws.Range("B25").Value = '"Website ID 56"+"\"+"3379\"'
$ws.Range("C25").Value = '"Website ID 56"+"\"+"3379\"'

$ws.Range("B26").Value = '"Website ID 56"+"\"+"SSA-5665\"'
$ws.Range("C26").Value = '"Website ID 56"+"\"+"SSA-5665\"'

$ws.Range("B29").Value = '"Website ID 56"+"\"+"Compassionate Allowances\"'
$ws.Range("C29").Value = '"Website ID 56"+"\"+"Compassionate Allowances\"'

# --- Rows 17-19, 27-28: column B then column C (same value) ---
$ws.Range("B17").Value = '"Website ID 56"+"\"+"SSA-4814\"'
$ws.Range("C17").Value = '"Website ID 56"+"\"+"SSA-4814\"'

$ws.Range("B18").Value = '"Website ID 56"+"\"+"Listing of Impairments\"'
$ws.Range("C18").Value = '"Website ID 56"+"\"+"Listing of Impairments\"'

$ws.Range("B19").Value = '"Website ID 56"+"\"+"Online Disability Appeal Application\"'
$ws.Range("C19").Value = '"Website ID 56"+"\"+"Online Disability Appeal Application\"'

$ws.Range("B27").Value = '"Website ID 56"+"\"+"Listing of Impairments\"'
$ws.Range("C27").Value = '"Website ID 56"+"\"+"Listing of Impairments\"'

$ws.Range("B28").Value = '"Website ID 56"+"\"+"Medicaid Eligibility Income\"'
$ws.Range("C28").Value = '"Website ID 56"+"\"+"Medicaid Eligibility Income\"'

# --- Row 2 (special typo in column C) and column C header ---
$ws.Range("C2").Value = '"Website ID 56+"\"+"SSA-3288\"'
$ws.Range("C1").Value = "Folder Path 2"

# --- Column B for rows 2-6 (unquoted style, filled last) ---
$ws.Range("B2").Value = "Website ID 56+\+SSA-3288\"
$ws.Range("B3").Value = "Website ID 56+\+SSA-827\"
$ws.Range("B4").Value = "Website ID 56+\+SSA-1696\"
$ws.Range("B5").Value = "Website ID 56+\+SSA-8000\"
$ws.Range("B6").Value = "Website ID 56+\+Fillable SSA-8000\"

# --- Column widths (closest achievable values to the target 57.21875 / 53.33203125) ---
$ws.Columns.Item(2).ColumnWidth = 56.25
$ws.Columns.Item(3).ColumnWidth = 52.5

# --- Selection as in the edited workbook ---
$ws.Range("B9").Select()
